$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username value in I2 to the new test user string
$ws.Range("I2").Value = "tcbdemotestuser2"

# Adjust column I width to fit the new, longer text (bestFit-style autofit).
# Compensate for the fixed ~0.8333 "cell padding" that Excel/this engine
# adds on top of whatever ColumnWidth is assigned, so the persisted width
# in the workbook ends up matching the target as closely as possible.
$ws.Columns.Item(9).ColumnWidth = (17.42578125 - 0.8333333333333334)

# Update the active selection to I2 (was J4)
$ws.Range("I2").Select() | Out-Null
